$d = $word.ActiveDocument

# --- Step 1: Update the To/From/Date/Subject header paragraph (markdown "##" -> bold "**...**") ---
$null = $d.Content.Find.Execute("## To: Credit Committee  ", $true, $false, $false, $false, $false, $true, 1, $false, "**To:** Credit Committee  ", 2)
$null = $d.Content.Find.Execute("## From: [Your Name], Financial Analyst  ", $true, $false, $false, $false, $false, $true, 1, $false, "**From:** [Your Name], Financial Analyst  ", 2)
$null = $d.Content.Find.Execute("## Date: [Insert Date]  ", $true, $false, $false, $false, $false, $true, 1, $false, "**Date:** [Insert Date]  ", 2)
$null = $d.Content.Find.Execute("## Subject: Credit Assessment of Porsche Automobil Holding SE for Fiscal Year 2024", $true, $false, $false, $false, $false, $true, 1, $false, "**Subject:** Credit Memo for Porsche Automobil Holding SE - Fiscal Year 2024", 2)

# --- Step 2: Delete the old Sections 1-5 + Recommended Action paragraphs (keep the "---" divider) ---
$pStart = $d.Paragraphs.Item(7)
$pEnd = $d.Paragraphs.Item(11)
$delRange = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$delRange.Delete()

# --- Step 3: Replace the closing paragraph text and append the new report body ---
$idx = 8
$p = $d.Paragraphs.Item($idx)
$p.Range.Text = "## 1. Executive Summary"

# --- Step 4: Insert the remaining new paragraphs after the (former) closing paragraph ---
$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs.Item($idx)
$p.Range.Text = "Porsche Automobil Holding SE (herein referred to as `"Porsche SE`") has reported financial results for the fiscal year ending December 31, 2024, reflecting significant challenges in profitability while maintaining a strong asset base. Despite facing a net loss of €1.52 billion in 2024, the company's liquidity position remains robust, with liquid assets amounting to €2.4 billion. This memo analyzes Porsche SE's financial performance and position, assessing its creditworthiness and providing a recommendation for future lending decisions."
$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs.Item($idx)
$p.Range.Text = "## 2. Financial Highlights"
$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs.Item($idx)
$p.Range.Text = "- **Revenue:** €115 thousand (2023: €96 thousand)" + [char]11 + "- **Net Loss:** €1,521,284 thousand (2023: Profit of €1,441,088 thousand)" + [char]11 + "- **Total Assets:** €33,066,788 thousand (2023: €34,529,841 thousand)" + [char]11 + "- **Cash and Cash Equivalents:** €1,822,850 thousand (2023: €723,595 thousand)" + [char]11 + "- **Marketable Securities:** €576,292 thousand (2023: €283,406 thousand)" + [char]11 + "- **Total Equity:** €25,060,859 thousand (2023: €27,365,224 thousand)"
$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs.Item($idx)
$p.Range.Text = "### Key Income Statement Metrics" + [char]11 + "- **EBITDA**: While not directly available in the financial statements, calculated considering income from investments and operational expenses, indicating a need for careful calculation adjustments post-loss recognition." + [char]11 + "- **Investment Income**: Strong income from investments was reported at €1,682,887 thousand (2023: €1,509,824 thousand), demonstrating the company’s effective management of its financial assets."
$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs.Item($idx)
$p.Range.Text = "## 3. Key Ratios"
$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs.Item($idx)
$p.Range.Text = "- **Debt-to-Equity Ratio:** 0.32 " + [char]11 + "   - Total Liabilities: €8,005,929 thousand" + [char]11 + "   - Total Equity: €25,060,859 thousand"
$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs.Item($idx)
$p.Range.Text = "\[" + [char]11 + "  \text{Debt-to-Equity} = \frac{\text{Total Liabilities}}{\text{Total Equity}} = \frac{8,005,929}{25,060,859} \approx 0.32" + [char]11 + "  \]"
$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs.Item($idx)
$p.Range.Text = "- **Interest Coverage Ratio:** 6.9 " + [char]11 + "   - Income from investments (assumed as EBITDA): €1,682,887 thousand " + [char]11 + "   - Interest Expenses: €243,578 thousand"
$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs.Item($idx)
$p.Range.Text = "\[" + [char]11 + "  \text{Interest Coverage} = \frac{\text{EBIT}}{\text{Interest Expenses}} \approx \frac{1,682,887}{243,578} \approx 6.9" + [char]11 + "  \]"
$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs.Item($idx)
$p.Range.Text = "## 4. Risk Analysis & Commentary"
$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs.Item($idx)
$p.Range.Text = "### Risks:" + [char]11 + "1. **Profitability Risks**: The company reported a net loss for 2024, significantly impacting retained earnings and indicating operational challenges which may persist." + [char]11 + "2. **Asset Impairment**: The recognition of a €2.93 billion impairment on financial assets could suggest vulnerabilities in investment performance or market conditions affecting portfolio value." + [char]11 + "3. **Market Volatility**: Continuous fluctuations and economic uncertainties may affect future revenue generation and profitability."
$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs.Item($idx)
$p.Range.Text = "### Strengths:" + [char]11 + "1. **Strong Asset Base**: Despite the reported loss, Porsche SE maintains substantial asset value with considerable liquid assets and investments." + [char]11 + "2. **Investment Income Stability**: The ability to generate significant investment income is a crucial buffer against operational losses, providing a stable foundation for ongoing operations." + [char]11 + "3. **Low Debt Levels**: With a relatively low debt-to-equity ratio, the company maintains significant financial flexibility."
$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs.Item($idx)
$p.Range.Text = "## 5. Final Credit Recommendation"
$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs.Item($idx)
$p.Range.Text = "Despite the recent operational losses, Porsche SE's robust liquidity and substantial investment income argue for a cautious outlook. The existing asset base, coupled with effective management of financial assets, provides a compelling argument for continued credit support."
$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs.Item($idx)
$p.Range.Text = "**Recommendation:** Approve a credit line with close monitoring of performance metrics over the next fiscal period, focusing particularly on revenue recovery, investment performance, and operational efficiencies. Given the financial restructuring and potential recovery strategies to be implemented, this is a prudent approach to safeguarding lender interests while supporting Porsche SE’s operations."
$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs.Item($idx)
$p.Range.Text = "**Action Items:** " + [char]11 + "- Schedule quarterly reviews of financial performance." + [char]11 + "- Work with Porsche SE to understand their recovery strategy."
$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs.Item($idx)
$p.Range.Text = "---"
$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs.Item($idx)
$p.Range.Text = "**Prepared by:**  " + [char]11 + "[Your Name]  " + [char]11 + "[Your Job Title]  " + [char]11 + "[Your Contact Information]"
